# Applies the "change logic to decide end game and winner" edit:
#  1. Strike-through the "when alive anti > alive norm, antis win" bullet.
#  2. Insert a brand-new bullet ("when there is only 1 norm left, ...")
#     before the "when alive users = 2" bullet, and strike-through +
#     extend the "when alive users = 2" bullet.
#  3. Strike-through the "if there is blank alive, that blank wins" bullet.
#  4. Strike-through the "if there is anti alive, antis win" bullet
#     (keeps its spell-check proofErr wrapper runs intact).
#  5. Split "If only Civilians are left in the game, Civilians win." into
#     three runs with identical text.
#  6. Replace the "Undercovers win." / "2 players left" bullets with a
#     single new bullet made of many small runs.

$d = $word.ActiveDocument

$xmlHeader = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$xmlFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

function Get-ParaIndexByText([string]$text) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $t = $d.Paragraphs.Item($i).Range.Text
        if ($t.TrimEnd() -eq $text) {
            return $i
        }
    }
    throw "Could not find paragraph with text: $text"
}

function Set-ParaStrike([int]$idx) {
    $p = $d.Paragraphs.Item($idx)
    $p.Range.Font.StrikeThrough = 1
}

# Replaces the paragraph(s) spanning from $startIdx to $endIdx (inclusive,
# 1-based Paragraphs indices) with the supplied raw <w:p>...</w:p> XML
# (one or more paragraphs back to back). Returns nothing; paragraph
# indices after this call must be re-resolved by text lookup since the
# paragraph count may change.
function Replace-Paragraphs([int]$startIdx, [int]$endIdx, [string]$paragraphsXml) {
    $pStart = $d.Paragraphs.Item($startIdx)
    $pEnd = $d.Paragraphs.Item($endIdx)
    $target = $d.Range($pStart.Range.Start, $pEnd.Range.End)
    $xml = $xmlHeader + $paragraphsXml + $xmlFooter
    $target.InsertXML($xml)
}

# --- 1. "when alive anti > alive norm, antis win" -> strike-through ---
$idx = Get-ParaIndexByText("when alive anti > alive norm, antis win")
Set-ParaStrike($idx)

# --- 2. insert new bullet + extend/strike "when alive users = 2" ---
$idx = Get-ParaIndexByText("when alive users = 2")
$newParas = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>when there is only 1 norm left,</w:t></w:r><w:r><w:t xml:space="preserve"> the game ends</w:t></w:r><w:r><w:t>. If there is blank alive, the blanks win. If not,</w:t></w:r><w:r><w:t xml:space="preserve"> the antis</w:t></w:r><w:r><w:t xml:space="preserve"> win.</w:t></w:r></w:p>' + `
    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0"/><w:rPr><w:strike/></w:rPr></w:pPr><w:r><w:rPr><w:strike/></w:rPr><w:t>when alive users = 2</w:t></w:r><w:r><w:rPr><w:strike/></w:rPr><w:t>, game ends and the alive blank/anti wins</w:t></w:r></w:p>'
Replace-Paragraphs $idx $idx $newParas

# --- 3. "if there is blank alive, that blank wins" -> strike-through ---
$idx = Get-ParaIndexByText("if there is blank alive, that blank wins")
Set-ParaStrike($idx)

# --- 4. "if there is anti alive, antis win" -> strike-through (text unchanged) ---
$idx = Get-ParaIndexByText("if there is anti alive, antis win")
Set-ParaStrike($idx)

# --- 5. "If only Civilians are left in the game, Civilians win." -> 3 runs ---
$idx = Get-ParaIndexByText("If only Civilians are left in the game, Civilians win.")
$newPara = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">If only Civilians are left in the game, </w:t></w:r><w:r><w:t xml:space="preserve">the </w:t></w:r><w:r><w:t>Civilians win.</w:t></w:r></w:p>'
Replace-Paragraphs $idx $idx $newPara

# --- 6. replace the two Undercover-related bullets with one combined bullet ---
$idx1 = Get-ParaIndexByText("If the number of Undercovers is more than that of Civilians, Undercovers win.")
$idx2 = Get-ParaIndexByText("If there only 2 players left in the game, if there is a Blank left, that Blank wins. If there is an Undercover left, Undercovers win.")
$newPara = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="3"/></w:numPr></w:pPr>' + `
    '<w:r><w:t xml:space="preserve">If only 1 Civilian is left, the </w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">game ends. If </w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">there </w:t></w:r>' + `
    '<w:r><w:t>are B</w:t></w:r>' + `
    '<w:r><w:t>lank</w:t></w:r>' + `
    '<w:r><w:t>s</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> alive, the </w:t></w:r>' + `
    '<w:r><w:t>Blanks</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> win. If not, the </w:t></w:r>' + `
    '<w:r><w:t>Undercovers</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> win.</w:t></w:r>' + `
    '</w:p>'
Replace-Paragraphs $idx1 $idx2 $newPara

Write-Output "Edit applied."
